$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-keyed updates: D = Price column, E = Volume(1h) column.
# Row 1 is the header; row 45 (USDe) is unchanged in this refresh and is omitted.
$updates = @(
    @{ Row = 2; D = "60.111.17"; E = "  -4.96%  " }
    @{ Row = 3; D = "2.996.03"; E = "  -5.34%  " }
    @{ Row = 4; D = $null; E = "  +0.15%  " }
    @{ Row = 5; D = "570.11"; E = "  -4.88%  " }
    @{ Row = 6; D = "124.99"; E = "  -7.51%  " }
    @{ Row = 7; D = $null; E = "  +0.23%  " }
    @{ Row = 8; D = "2.988.28"; E = "  -5.56%  " }
    @{ Row = 9; D = $null; E = "  -2.13%  " }
    @{ Row = 10; D = $null; E = "  -7.38%  " }
    @{ Row = 11; D = "5.04"; E = "  -5.13%  " }
    @{ Row = 12; D = $null; E = "  -2.63%  " }
    @{ Row = 13; D = $null; E = "  -7.58%  " }
    @{ Row = 14; D = "32.46"; E = "  -6.07%  " }
    @{ Row = 15; D = $null; E = "  +0.41%  " }
    @{ Row = 16; D = "3.490.67"; E = "  -5.30%  " }
    @{ Row = 17; D = "2.998.31"; E = "  -5.17%  " }
    @{ Row = 18; D = "60.107.56"; E = "  -4.94%  " }
    @{ Row = 19; D = "6.52"; E = "  -0.55%  " }
    @{ Row = 20; D = "429.01"; E = "  -6.72%  " }
    @{ Row = 21; D = "13.12"; E = "  -5.68%  " }
    @{ Row = 22; D = "0.671"; E = "  -3.30%  " }
    @{ Row = 23; D = $null; E = "  -7.39%  " }
    @{ Row = 24; D = $null; E = "  -2.27%  " }
    @{ Row = 25; D = "79.48"; E = "  -4.23%  " }
    @{ Row = 26; D = $null; E = "  +0.10%  " }
    @{ Row = 27; D = $null; E = "  +0.14%  " }
    @{ Row = 28; D = $null; E = "  -6.32%  " }
    @{ Row = 29; D = $null; E = "  -4.94%  " }
    @{ Row = 30; D = "7.18"; E = "  -6.78%  " }
    @{ Row = 31; D = "6.10"; E = "  -10.08%  " }
    @{ Row = 32; D = "25.24"; E = "  -7.09%  " }
    @{ Row = 33; D = "0.0951"; E = "  -5.31%  " }
    @{ Row = 34; D = $null; E = "  -4.66%  " }
    @{ Row = 35; D = $null; E = "  -8.71%  " }
    @{ Row = 36; D = "50.24"; E = "  -2.06%  " }
    @{ Row = 37; D = $null; E = "  -16.16%  " }
    @{ Row = 38; D = "8.49"; E = "  +4.47%  " }
    @{ Row = 39; D = $null; E = "  -10.06%  " }
    @{ Row = 40; D = "0.0355"; E = "  -8.70%  " }
    @{ Row = 41; D = $null; E = "  -4.40%  " }
    @{ Row = 42; D = "371.28"; E = "  -5.28%  " }
    @{ Row = 43; D = "2.668.35"; E = "  -4.30%  " }
    @{ Row = 44; D = $null; E = "  -6.55%  " }
    @{ Row = 46; D = "121.43"; E = "  -4.45%  " }
    @{ Row = 47; D = "0.233"; E = "  -6.99%  " }
    @{ Row = 48; D = $null; E = "  -6.09%  " }
    @{ Row = 49; D = $null; E = "  -3.50%  " }
    @{ Row = 50; D = "23.26"; E = "  -6.81%  " }
    @{ Row = 51; D = "0.132"; E = "  -2.82%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        # Preserve these Price values as text (they use "." as a thousands
        # separator, e.g. "60.111.17"), matching the source data's inline-string type.
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

Write-Output "Updated $($updates.Count) rows"
